# Auto-generated edit script applying the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $Addr, $Text) {
    $cell = $Sheet.Range($Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "59.150.27"
Set-TextValue $ws "E2" "  -0.51%  "
Set-TextValue $ws "D3" "2.526.86"
Set-TextValue $ws "E3" "  +0.48%  "
Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  -0.04%  "
Set-TextValue $ws "D5" "537.50"
Set-TextValue $ws "E5" "  -0.63%  "
Set-TextValue $ws "D6" "137.22"
Set-TextValue $ws "E6" "  -1.75%  "
Set-TextValue $ws "D7" "0.999"
Set-TextValue $ws "E7" "  +0.08%  "
Set-TextValue $ws "D8" "0.571"
Set-TextValue $ws "E8" "  +1.12%  "
Set-TextValue $ws "D9" "2.526.78"
Set-TextValue $ws "E9" "  +0.39%  "
Set-TextValue $ws "E10" "  -0.27%  "
Set-TextValue $ws "D11" "0.158"
Set-TextValue $ws "E11" "  -1.85%  "
Set-TextValue $ws "D12" "5.34"
Set-TextValue $ws "E12" "  -0.63%  "
Set-TextValue $ws "E13" "  -2.92%  "
Set-TextValue $ws "D14" "2.968.02"
Set-TextValue $ws "E14" "  +0.13%  "
Set-TextValue $ws "D15" "23.09"
Set-TextValue $ws "E15" "  -1.23%  "
Set-TextValue $ws "D16" "58.852.64"
Set-TextValue $ws "E16" "  -0.80%  "
Set-TextValue $ws "E17" "  -1.36%  "
Set-TextValue $ws "D18" "2.514.63"
Set-TextValue $ws "E18" "  +0.02%  "
Set-TextValue $ws "D19" "11.17"
Set-TextValue $ws "E19" "  +0.64%  "
Set-TextValue $ws "E20" "  -0.35%  "
Set-TextValue $ws "D21" "323.85"
Set-TextValue $ws "E21" "  -0.47%  "
Set-TextValue $ws "E22" "  +0.02%  "
Set-TextValue $ws "E23" "  +1.97%  "
Set-TextValue $ws "D24" "65.83"
Set-TextValue $ws "E24" "  +3.98%  "
Set-TextValue $ws "E25" "  +0.17%  "
Set-TextValue $ws "E26" "  -1.86%  "
Set-TextValue $ws "E27" "  -0.36%  "
Set-TextValue $ws "E28" "  -3.44%  "
Set-TextValue $ws "B29" "PEPE"
Set-TextValue $ws "C29" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws "D29" "0.0₃0775"
Set-TextValue $ws "E29" "  -0.71%  "
Set-TextValue $ws "B30" "Aptos"
Set-TextValue $ws "C30" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws "D30" "6.71"
Set-TextValue $ws "E30" "  -1.95%  "
Set-TextValue $ws "E31" "  -1.28%  "
Set-TextValue $ws "D32" "168.03"
Set-TextValue $ws "E32" "  +2.41%  "
Set-TextValue $ws "E33" "  +5.36%  "
Set-TextValue $ws "D34" "0.999"
Set-TextValue $ws "E34" "  +0.01%  "
Set-TextValue $ws "E35" "  +1.97%  "
Set-TextValue $ws "D36" "18.45"
Set-TextValue $ws "E36" "  -0.31%  "
Set-TextValue $ws "D37" "4.13"
Set-TextValue $ws "E37" "  -2.82%  "
Set-TextValue $ws "D38" "1.55"
Set-TextValue $ws "E38" "  -3.15%  "
Set-TextValue $ws "E39" "  -0.66%  "
Set-TextValue $ws "D40" "0.815"
Set-TextValue $ws "E40" "  +0.46%  "
Set-TextValue $ws "D41" "3.63"
Set-TextValue $ws "E41" "  -1.42%  "
Set-TextValue $ws "D42" "285.28"
Set-TextValue $ws "E42" "  +1.15%  "
Set-TextValue $ws "D43" "5.17"
Set-TextValue $ws "E43" "  -1.44%  "
Set-TextValue $ws "D44" "132.86"
Set-TextValue $ws "E44" "  +6.51%  "
Set-TextValue $ws "E45" "  +0.20%  "
Set-TextValue $ws "E46" "  +1.56%  "
Set-TextValue $ws "D47" "10.87"
Set-TextValue $ws "E47" "  -0.01%  "
Set-TextValue $ws "E48" "  -1.15%  "
Set-TextValue $ws "D49" "0.0509"
Set-TextValue $ws "E49" "  -0.99%  "
Set-TextValue $ws "E50" "  -2.03%  "
Set-TextValue $ws "D51" "17.41"
Set-TextValue $ws "E51" "  -2.80%  "
